# Insert a new "mu" column before the existing "isSelected" column (column I).
# This shifts isSelected, bandwidth, transRate, uploadTime, totalTime one
# column to the right (I->J, J->K, K->L, L->M, M->N) and the new column I
# is populated with the same values as the "dataSize" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at I, pushing isSelected..totalTime to J..N
$ws.Columns.Item(9).Insert()

# Header for the new column
$ws.Cells.Item(1, 9).Value = "mu"

# Fill the new "mu" column (I) with the same values as "dataSize" (D)
$lastRow = $ws.Cells.Item(1, 4).End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 4).Value2
}
